$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Fitness) values change in blocks based on Generation (column B / row number)
# Rows 2-6   -> 7647
# Rows 7-35  -> 7320
# Rows 36-50 -> 7295
# Rows 51-252 -> 7293

$ws.Range("C2:C6").Value = 7647
$ws.Range("C7:C35").Value = 7320
$ws.Range("C36:C50").Value = 7295
$ws.Range("C51:C252").Value = 7293
